# Refresh cryptocurrency price/volume snapshot (GitHub Actions scheduled update).
# Column D ("Price") and column E ("Volume(1h)") are stored as plain text in the
# sheet, so numeric-looking price strings are written with a leading apostrophe
# (PowerShell/Excel text-prefix convention) to stop Excel from auto-converting them
# into numbers, matching the original inlineStr cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.109.92'
$ws.Range('E2').Value = '  -1.21%  '

$ws.Range('D3').Value = '1.793.82'
$ws.Range('E3').Value = '  -0.59%  '

$ws.Range('D4').Value = "'1.001"
$ws.Range('E4').Value = '  +0.13%  '

$ws.Range('D5').Value = "'316.93"
$ws.Range('E5').Value = '  +0.55%  '

$ws.Range('D6').Value = "'1.000"
$ws.Range('E6').Value = '  +0.08%  '

$ws.Range('D7').Value = "'0.5363"
$ws.Range('E7').Value = '  -2.55%  '

$ws.Range('D8').Value = "'0.3769"
$ws.Range('E8').Value = '  -2.29%  '

$ws.Range('D9').Value = "'0.07471"

$ws.Range('D10').Value = "'41.70"
$ws.Range('E10').Value = '  -2.06%  '

$ws.Range('D11').Value = "'1.094"
$ws.Range('E11').Value = '  -3.01%  '

$ws.Range('D12').Value = "'1.001"
$ws.Range('E12').Value = '  +0.16%  '

$ws.Range('D13').Value = "'20.62"
$ws.Range('E13').Value = '  -2.93%  '

$ws.Range('D14').Value = "'6.106"
$ws.Range('E14').Value = '  -1.51%  '

$ws.Range('D15').Value = "'7.211"
$ws.Range('E15').Value = '  -3.49%  '

$ws.Range('D16').Value = '1.778.23'
$ws.Range('E16').Value = '  -1.60%  '

$ws.Range('D17').Value = "'89.15"
$ws.Range('E17').Value = '  -3.19%  '

$ws.Range('D18').Value = "'0.00001057"
$ws.Range('E18').Value = '  -1.45%  '

$ws.Range('D19').Value = "'0.06455"
$ws.Range('E19').Value = '  +0.19%  '

$ws.Range('D20').Value = "'0.9997"
$ws.Range('E20').Value = '  +0.02%  '

$ws.Range('D21').Value = "'17.34"
$ws.Range('E21').Value = '  -0.13%  '

$ws.Range('E22').Value = '  -1.22%  '

$ws.Range('D23').Value = '28.136.11'
$ws.Range('E23').Value = '  -1.13%  '

$ws.Range('E24').Value = '  -2.14%  '

$ws.Range('D25').Value = "'2.095"
$ws.Range('E25').Value = '  -2.10%  '

$ws.Range('D26').Value = "'154.78"
$ws.Range('E26').Value = '  -2.67%  '

$ws.Range('E27').Value = '  -2.40%  '

$ws.Range('D28').Value = '1.991.27'
$ws.Range('E28').Value = '  -1.13%  '

$ws.Range('D29').Value = "'2.273"
$ws.Range('E29').Value = '  -5.85%  '

$ws.Range('D30').Value = "'120.64"
$ws.Range('E30').Value = '  -2.69%  '

$ws.Range('E31').Value = '  -0.75%  '

$ws.Range('D32').Value = "'0.1055"
$ws.Range('E32').Value = '  +3.13%  '

$ws.Range('D33').Value = "'3.654"
$ws.Range('E33').Value = '  -0.90%  '

$ws.Range('D34').Value = "'5.561"
$ws.Range('E34').Value = '  -3.94%  '

$ws.Range('D35').Value = "'0.06547"
$ws.Range('E35').Value = '  +1.38%  '

$ws.Range('D36').Value = "'0.2257"
$ws.Range('E36').Value = '  -2.60%  '

$ws.Range('D37').Value = "'0.02283"
$ws.Range('E37').Value = '  -1.94%  '

$ws.Range('D38').Value = "'5.027"
$ws.Range('E38').Value = '  -2.99%  '

$ws.Range('D39').Value = "'8.452"
$ws.Range('E39').Value = '  -4.30%  '

$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = "'0.6165"
$ws.Range('E40').Value = '  -3.93%  '

$ws.Range('B41').Value = 'WEMIXTOKEN'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').Value = "'1.447"
$ws.Range('E41').Value = '  +4.50%  '

$ws.Range('D42').Value = "'11.08"
$ws.Range('E42').Value = '  -5.11%  '

$ws.Range('E43').Value = '  +0.97%  '

$ws.Range('D44').Value = "'0.9994"
$ws.Range('E44').Value = '  +0.03%  '

$ws.Range('D45').Value = "'13.25"
$ws.Range('E45').Value = '  -2.21%  '

$ws.Range('E46').Value = '  -0.21%  '

$ws.Range('D47').Value = "'0.5781"
$ws.Range('E47').Value = '  -3.51%  '

$ws.Range('D48').Value = "'127.50"
$ws.Range('E48').Value = '  +0.20%  '

$ws.Range('D49').Value = "'1.189"
$ws.Range('E49').Value = '  +3.25%  '

$ws.Range('D50').Value = "'1.927"
$ws.Range('E50').Value = '  -2.89%  '

$ws.Range('D51').Value = "'0.06813"
$ws.Range('E51').Value = '  -1.17%  '
